$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================
# Fix fixture ordering / odds data for Italy Serie C Group C
# 2023-2024 (rows 15-18, 75-77, 90, 92-95) and append the
# three newly scraped matches (rows 97-99).
# =========================================================

# --- Update existing rows (re-synced fixture/odds data) ---
# Row 15
$ws.Range("F15").Value = "Benevento"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = "Virtus Francavilla"
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1.56
$ws.Range("K15").Value = "07/09/2023 15:12"
$ws.Range("L15").Value = 1.56
$ws.Range("M15").Value = "11/09/2023 20:36"
$ws.Range("N15").Value = 3.69
$ws.Range("O15").Value = "07/09/2023 15:12"
$ws.Range("P15").Value = 3.84
$ws.Range("Q15").Value = "11/09/2023 20:36"
$ws.Range("R15").Value = 5.31
$ws.Range("S15").Value = "07/09/2023 15:12"
$ws.Range("T15").Value = 6.45
$ws.Range("U15").Value = "11/09/2023 20:37"
$ws.Range("V15").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/benevento-virtus-francavilla/QZ3JJmpl/"

# Row 16
$ws.Range("F16").Value = "Picerno"
$ws.Range("G16").Value = 1
$ws.Range("H16").Value = "Taranto"
$ws.Range("I16").Value = 1
$ws.Range("J16").Value = 1.96
$ws.Range("K16").Value = "07/09/2023 15:12"
$ws.Range("L16").Value = 2.4
$ws.Range("M16").Value = "11/09/2023 20:43"
$ws.Range("N16").Value = 2.95
$ws.Range("O16").Value = "07/09/2023 15:12"
$ws.Range("P16").Value = 2.56
$ws.Range("Q16").Value = "11/09/2023 20:36"
$ws.Range("R16").Value = 3.96
$ws.Range("S16").Value = "07/09/2023 15:12"
$ws.Range("T16").Value = 4
$ws.Range("U16").Value = "11/09/2023 20:43"
$ws.Range("V16").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/picerno-taranto/AyqOMyIP/"

# Row 17
$ws.Range("F17").Value = "Juve Stabia"
$ws.Range("G17").Value = 1
$ws.Range("H17").Value = "Avellino"
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2.49
$ws.Range("K17").Value = "07/09/2023 15:12"
$ws.Range("L17").Value = 2.84
$ws.Range("M17").Value = "11/09/2023 20:22"
$ws.Range("N17").Value = 2.85
$ws.Range("O17").Value = "07/09/2023 15:12"
$ws.Range("P17").Value = 2.91
$ws.Range("Q17").Value = "11/09/2023 19:09"
$ws.Range("R17").Value = 2.96
$ws.Range("S17").Value = "07/09/2023 15:12"
$ws.Range("T17").Value = 2.77
$ws.Range("U17").Value = "11/09/2023 20:22"
$ws.Range("V17").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/juve-stabia-avellino/UR0zG5qD/"

# Row 18
$ws.Range("F18").Value = "Foggia"
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = "Giugliano"
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 1.8
$ws.Range("K18").Value = "07/09/2023 15:12"
$ws.Range("L18").Value = 2.05
$ws.Range("M18").Value = "11/09/2023 20:35"
$ws.Range("N18").Value = 3.29
$ws.Range("O18").Value = "07/09/2023 15:12"
$ws.Range("P18").Value = 3.43
$ws.Range("Q18").Value = "11/09/2023 20:35"
$ws.Range("R18").Value = 4.3
$ws.Range("S18").Value = "07/09/2023 15:12"
$ws.Range("T18").Value = 3.58
$ws.Range("U18").Value = "11/09/2023 20:35"
$ws.Range("V18").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/foggia-giugliano/0hdWGoU6/"

# Row 75
$ws.Range("F75").Value = "Latina"
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = "Monopoli"
$ws.Range("I75").Value = 2
$ws.Range("J75").Value = 1.98
$ws.Range("K75").Value = "12/10/2023 08:13"
$ws.Range("L75").Value = 1.93
$ws.Range("M75").Value = "15/10/2023 16:11"
$ws.Range("N75").Value = 3
$ws.Range("O75").Value = "12/10/2023 08:13"
$ws.Range("P75").Value = 3.13
$ws.Range("Q75").Value = "15/10/2023 16:11"
$ws.Range("R75").Value = 3.79
$ws.Range("S75").Value = "12/10/2023 08:13"
$ws.Range("T75").Value = 4.52
$ws.Range("U75").Value = "15/10/2023 16:11"
$ws.Range("V75").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/latina-monopoli/0lxtjXss/"

# Row 76
$ws.Range("F76").Value = "Taranto"
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = "Crotone"
$ws.Range("I76").Value = 1
$ws.Range("J76").Value = 3.02
$ws.Range("K76").Value = "12/10/2023 08:13"
$ws.Range("L76").Value = 3.66
$ws.Range("M76").Value = "15/10/2023 16:09"
$ws.Range("N76").Value = 2.87
$ws.Range("O76").Value = "12/10/2023 08:13"
$ws.Range("P76").Value = 2.71
$ws.Range("Q76").Value = "15/10/2023 16:09"
$ws.Range("R76").Value = 2.38
$ws.Range("S76").Value = "12/10/2023 08:13"
$ws.Range("T76").Value = 2.41
$ws.Range("U76").Value = "15/10/2023 16:09"
$ws.Range("V76").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/taranto-crotone/WngimVRa/"

# Row 77
$ws.Range("F77").Value = "Virtus Francavilla"
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = "Sorrento"
$ws.Range("I77").Value = 5
$ws.Range("J77").Value = 1.71
$ws.Range("K77").Value = "12/10/2023 08:13"
$ws.Range("L77").Value = 2.17
$ws.Range("M77").Value = "15/10/2023 16:10"
$ws.Range("N77").Value = 3.32
$ws.Range("O77").Value = "12/10/2023 08:13"
$ws.Range("P77").Value = 2.97
$ws.Range("Q77").Value = "15/10/2023 16:10"
$ws.Range("R77").Value = 4.61
$ws.Range("S77").Value = "12/10/2023 08:13"
$ws.Range("T77").Value = 3.85
$ws.Range("U77").Value = "15/10/2023 16:10"
$ws.Range("V77").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/virtus-francavilla-sorrento/EBrdnks6/"

# Row 90
$ws.Range("F90").Value = "Avellino"
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = "Audace Cerignola"
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 1.76
$ws.Range("K90").Value = "24/10/2023 12:42"
$ws.Range("L90").Value = 2.01
$ws.Range("M90").Value = "25/10/2023 18:29"
$ws.Range("N90").Value = 3.18
$ws.Range("O90").Value = "24/10/2023 12:42"
$ws.Range("P90").Value = 3.07
$ws.Range("Q90").Value = "25/10/2023 18:29"
$ws.Range("R90").Value = 4.52
$ws.Range("S90").Value = "24/10/2023 12:42"
$ws.Range("T90").Value = 4.28
$ws.Range("U90").Value = "25/10/2023 18:29"
$ws.Range("V90").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/avellino-audace-cerignola/vZXZcntI/"

# Row 92
$ws.Range("F92").Value = "Latina"
$ws.Range("G92").Value = 1
$ws.Range("H92").Value = "Virtus Francavilla"
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 1.94
$ws.Range("K92").Value = "24/10/2023 12:42"
$ws.Range("L92").Value = 2.01
$ws.Range("M92").Value = "25/10/2023 18:28"
$ws.Range("N92").Value = 3.14
$ws.Range("O92").Value = "24/10/2023 12:42"
$ws.Range("P92").Value = 3.15
$ws.Range("Q92").Value = "25/10/2023 18:29"
$ws.Range("R92").Value = 3.87
$ws.Range("S92").Value = "24/10/2023 12:42"
$ws.Range("T92").Value = 4.12
$ws.Range("U92").Value = "25/10/2023 18:28"
$ws.Range("V92").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/latina-virtus-francavilla/CM2lH9tt/"

# Row 93
$ws.Range("F93").Value = "Monterosi"
$ws.Range("G93").Value = 1
$ws.Range("H93").Value = "Catania"
$ws.Range("I93").Value = 1
$ws.Range("J93").Value = 4.02
$ws.Range("K93").Value = "24/10/2023 12:42"
$ws.Range("L93").Value = 6.28
$ws.Range("M93").Value = "25/10/2023 18:27"
$ws.Range("N93").Value = 3.28
$ws.Range("O93").Value = "24/10/2023 12:42"
$ws.Range("P93").Value = 4
$ws.Range("Q93").Value = "25/10/2023 18:24"
$ws.Range("R93").Value = 1.83
$ws.Range("S93").Value = "24/10/2023 12:42"
$ws.Range("T93").Value = 1.54
$ws.Range("U93").Value = "25/10/2023 18:24"
$ws.Range("V93").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/monterosi-catania/YyicFmBh/"

# Row 94
$ws.Range("F94").Value = "Casertana"
$ws.Range("G94").Value = 2
$ws.Range("H94").Value = "Juve Stabia"
$ws.Range("I94").Value = 1
$ws.Range("J94").Value = 2.56
$ws.Range("K94").Value = "24/10/2023 12:42"
$ws.Range("L94").Value = 2.69
$ws.Range("M94").Value = "25/10/2023 20:37"
$ws.Range("N94").Value = 2.83
$ws.Range("O94").Value = "24/10/2023 12:42"
$ws.Range("P94").Value = 3.01
$ws.Range("Q94").Value = "25/10/2023 20:37"
$ws.Range("R94").Value = 2.8
$ws.Range("S94").Value = "24/10/2023 12:42"
$ws.Range("T94").Value = 2.83
$ws.Range("U94").Value = "25/10/2023 20:37"
$ws.Range("V94").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/casertana-juve-stabia/vL2bC8UL/"

# Row 95
$ws.Range("F95").Value = "Monopoli"
$ws.Range("G95").Value = 1
$ws.Range("H95").Value = "Picerno"
$ws.Range("I95").Value = 1
$ws.Range("J95").Value = 2.59
$ws.Range("K95").Value = "24/10/2023 12:42"
$ws.Range("L95").Value = 2.67
$ws.Range("M95").Value = "25/10/2023 20:41"
$ws.Range("N95").Value = 2.92
$ws.Range("O95").Value = "24/10/2023 12:42"
$ws.Range("P95").Value = 2.97
$ws.Range("Q95").Value = "25/10/2023 20:38"
$ws.Range("R95").Value = 2.69
$ws.Range("S95").Value = "24/10/2023 12:42"
$ws.Range("T95").Value = 2.9
$ws.Range("U95").Value = "25/10/2023 20:41"
$ws.Range("V95").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/monopoli-picerno/jD1hGTen/"

# --- Append new rows 97-99 (copy style from row 96 first) ---
$ws.Range("A96:V96").Copy() | Out-Null
$ws.Range("A97:V99").PasteSpecial(-4122) | Out-Null

# Row 97
$ws.Range("A97").Value = 96
$ws.Range("B97").Value = "italy"
$ws.Range("C97").Value = "serie-c-group-c"
$ws.Range("D97").Value = "2023-2024"
$ws.Range("E97").Value = 45225.77083333334
$ws.Range("F97").Value = "ACR Messina"
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = "Brindisi"
$ws.Range("I97").Value = 1
$ws.Range("J97").Value = 2.15
$ws.Range("K97").Value = "24/10/2023 21:12"
$ws.Range("L97").Value = 2.25
$ws.Range("M97").Value = "26/10/2023 18:28"
$ws.Range("N97").Value = 2.93
$ws.Range("O97").Value = "24/10/2023 21:12"
$ws.Range("P97").Value = 2.95
$ws.Range("Q97").Value = "26/10/2023 18:28"
$ws.Range("R97").Value = 3.38
$ws.Range("S97").Value = "24/10/2023 21:12"
$ws.Range("T97").Value = 3.65
$ws.Range("U97").Value = "26/10/2023 18:28"
$ws.Range("V97").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/acr-messina-brindisi/QVTVbSQB/"

# Row 98
$ws.Range("A98").Value = 97
$ws.Range("B98").Value = "italy"
$ws.Range("C98").Value = "serie-c-group-c"
$ws.Range("D98").Value = "2023-2024"
$ws.Range("E98").Value = 45225.86458333334
$ws.Range("F98").Value = "Foggia"
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = "Benevento"
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 2.27
$ws.Range("K98").Value = "24/10/2023 21:12"
$ws.Range("L98").Value = 2.93
$ws.Range("M98").Value = "26/10/2023 20:41"
$ws.Range("N98").Value = 2.87
$ws.Range("O98").Value = "24/10/2023 21:12"
$ws.Range("P98").Value = 2.89
$ws.Range("Q98").Value = "26/10/2023 20:39"
$ws.Range("R98").Value = 3.2
$ws.Range("S98").Value = "24/10/2023 21:12"
$ws.Range("T98").Value = 2.7
$ws.Range("U98").Value = "26/10/2023 20:41"
$ws.Range("V98").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/foggia-benevento/McNvc6eO/"

# Row 99
$ws.Range("A99").Value = 98
$ws.Range("B99").Value = "italy"
$ws.Range("C99").Value = "serie-c-group-c"
$ws.Range("D99").Value = "2023-2024"
$ws.Range("E99").Value = 45225.86458333334
$ws.Range("F99").Value = "Potenza"
$ws.Range("G99").Value = 3
$ws.Range("H99").Value = "Sorrento"
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1.88
$ws.Range("K99").Value = "24/10/2023 21:12"
$ws.Range("L99").Value = 1.95
$ws.Range("M99").Value = "26/10/2023 20:42"
$ws.Range("N99").Value = 3.05
$ws.Range("O99").Value = "24/10/2023 21:12"
$ws.Range("P99").Value = 3.2
$ws.Range("Q99").Value = "26/10/2023 20:43"
$ws.Range("R99").Value = 4.1
$ws.Range("S99").Value = "24/10/2023 21:12"
$ws.Range("T99").Value = 4.31
$ws.Range("U99").Value = "26/10/2023 20:43"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/italy/serie-c-group-c/potenza-sorrento/pCf1E7Qb/"
